$wb = $excel.ActiveWorkbook

# "Data" sheet: the browser used by the second scenario (row 2) switches
# from the Chrome grid to the Android driver.
$data = $wb.Worksheets.Item("Data")
$data.Range("C2").Value = "android"

# Appium/Android run now ends up being the one the user left selected:
# make "Data" the active sheet with D7 selected, and drop the selection
# highlight from "Test".
$data.Activate()
$data.Range("D7").Select()
